$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the raw metric values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 275124.77
$metrics.Range("B3").Value  = 225756.45
$metrics.Range("B4").Value  = 87627.419999999984
$metrics.Range("B5").Value  = 10951
$metrics.Range("B6").Value  = 4642256.24
$metrics.Range("B7").Value  = 3915575.1199999996
$metrics.Range("B8").Value  = 1358229.56
$metrics.Range("B9").Value  = 179952
$metrics.Range("B10").Value = 33107580.040999822
$metrics.Range("B11").Value = 31190796.640000001
$metrics.Range("B12").Value = 11639938.450000001
$metrics.Range("B13").Value = 1277579

# today sheet's B/E/F columns (formulas off Metrics!B2:B13) recalc automatically.

# --- Restore/update the saved selections on each sheet ---
$metrics.Activate() | Out-Null
$metrics.Range("F11").Select() | Out-Null

$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("I15").Select() | Out-Null
